$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-55 down to 40-56
$ws.Rows.Item(39).Insert()

# Copy style of date cell from the row that is now 40 (was row 39) to new row 39
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new row 39 with data (same market/category constants + new weekly record)
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44582
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = 100112031
$ws.Cells.Item(39, 7).Value = "Poroto verde"
$ws.Cells.Item(39, 8).Value = "Magnum"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 50
$ws.Cells.Item(39, 11).Value = 38000
$ws.Cells.Item(39, 12).Value = 38000
$ws.Cells.Item(39, 13).Value = 38000
$ws.Cells.Item(39, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Región Metropolitana"
$ws.Cells.Item(39, 16).Value = 1520
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
